# Auto-generated Excel COM-interop script
# Applies updated market-price snapshot values to the Coeurl_Profits workbook
# (columns H-N: currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ), LeveProfit(NQ/HQ))

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1205.2273
$ws.Range("I28").Value = 1132.8
$ws.Range("J28").Value = 1360.4286
$ws.Range("K28").Value = 1132.8
$ws.Range("L28").Value = 1360.4286
$ws.Range("M28").Value = -647.8
$ws.Range("N28").Value = -2330.4286
# Row 43
$ws.Range("H43").Value = 14429.75
$ws.Range("I43").Value = 50475
$ws.Range("K43").Value = 50475
$ws.Range("M43").Value = -50406
# Row 82
$ws.Range("H82").Value = 4285.25
$ws.Range("I82").Value = 4285.25
$ws.Range("K82").Value = 12855.75
$ws.Range("M82").Value = -12449.75
# Row 85
$ws.Range("H85").Value = 4285.25
$ws.Range("I85").Value = 4285.25
$ws.Range("K85").Value = 12855.75
$ws.Range("M85").Value = -11451.75
# Row 131
$ws.Range("H131").Value = 5561
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 10622
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 31866
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -41946
# Row 132
$ws.Range("H132").Value = 11564447
$ws.Range("I132").Value = 12383998
$ws.Range("J132").Value = 500502.5
$ws.Range("K132").Value = 37151994
$ws.Range("L132").Value = 1501507.5
$ws.Range("M132").Value = -37149464
$ws.Range("N132").Value = -1506567.5
# Row 137
$ws.Range("H137").Value = 5794.567
$ws.Range("I137").Value = 1290.6471
$ws.Range("J137").Value = 11684.308
$ws.Range("K137").Value = 3871.9413
$ws.Range("L137").Value = 35052.924
$ws.Range("M137").Value = -1321.9413
$ws.Range("N137").Value = -40152.924

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 22275.75
$ws.Range("I45").Value = 16379.5
$ws.Range("K45").Value = 16379.5
$ws.Range("M45").Value = -16002.5
# Row 74
$ws.Range("H74").Value = 9095.125
$ws.Range("I74").Value = 3037.2856
$ws.Range("K74").Value = 3037.2856
$ws.Range("M74").Value = -2163.2856
# Row 77
$ws.Range("H77").Value = 9095.125
$ws.Range("I77").Value = 3037.2856
$ws.Range("K77").Value = 15186.428
$ws.Range("M77").Value = -10818.428
# Row 102
$ws.Range("H102").Value = 3720.25
$ws.Range("I102").Value = 3556.0417
$ws.Range("J102").Value = 4048.6667
$ws.Range("K102").Value = 3556.0417
$ws.Range("L102").Value = 4048.6667
$ws.Range("M102").Value = -1934.0417
$ws.Range("N102").Value = -7292.6667
# Row 110
$ws.Range("H110").Value = 7760.28
$ws.Range("I110").Value = 10980.637
$ws.Range("K110").Value = 10980.637
$ws.Range("M110").Value = -8935.637000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 4028.6
$ws.Range("I20").Value = 3491.7144
$ws.Range("J20").Value = 6847.25
$ws.Range("K20").Value = 3491.7144
$ws.Range("L20").Value = 6847.25
$ws.Range("M20").Value = -3244.7144
$ws.Range("N20").Value = -7341.25
# Row 86
$ws.Range("H86").Value = 9117.647000000001
$ws.Range("I86").Value = 6300.2085
$ws.Range("K86").Value = 6300.2085
$ws.Range("M86").Value = -5177.2085
# Row 89
$ws.Range("H89").Value = 9117.647000000001
$ws.Range("I89").Value = 6300.2085
$ws.Range("K89").Value = 31501.0425
$ws.Range("M89").Value = -25885.0425
# Row 94
$ws.Range("H94").Value = 2222.7144
$ws.Range("I94").Value = 1829
$ws.Range("J94").Value = 3666.3333
$ws.Range("K94").Value = 1829
$ws.Range("L94").Value = 3666.3333
$ws.Range("M94").Value = -1378
$ws.Range("N94").Value = -4568.3333

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 3352.4285
$ws.Range("I16").Value = 3352.4285
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3352.4285
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3065.4285
$ws.Range("N16").Value = ""
# Row 31
$ws.Range("H31").Value = 11970039
$ws.Range("I31").Value = 23897938
$ws.Range("K31").Value = 23897938
$ws.Range("M31").Value = -23897643
# Row 34
$ws.Range("H34").Value = 11970039
$ws.Range("I34").Value = 23897938
$ws.Range("K34").Value = 23897938
$ws.Range("M34").Value = -23897736
# Row 99
$ws.Range("H99").Value = 3779.0952
$ws.Range("I99").Value = 3683.923
$ws.Range("K99").Value = 3683.923
$ws.Range("M99").Value = -2185.923
# Row 107
$ws.Range("H107").Value = 967.6667
$ws.Range("I107").Value = 1007.86664
$ws.Range("J107").Value = 766.6667
$ws.Range("K107").Value = 1007.86664
$ws.Range("L107").Value = 766.6667
$ws.Range("M107").Value = 912.13336
$ws.Range("N107").Value = -4606.6667
# Row 113
$ws.Range("H113").Value = 3352.4285
$ws.Range("I113").Value = 3352.4285
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3352.4285
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1182.4285
$ws.Range("N113").Value = ""
# Row 126
$ws.Range("H126").Value = 3779.0952
$ws.Range("I126").Value = 3683.923
$ws.Range("K126").Value = 11051.769
$ws.Range("M126").Value = -8581.769
# Row 133
$ws.Range("H133").Value = 47487.5
$ws.Range("I133").Value = 47450
$ws.Range("J133").Value = 47500
$ws.Range("K133").Value = 47450
$ws.Range("L133").Value = 47500
$ws.Range("M133").Value = -44920
$ws.Range("N133").Value = -52560
# Row 134
$ws.Range("H134").Value = 11207.429
$ws.Range("I134").Value = 2556.8696
$ws.Range("K134").Value = 7670.6088
$ws.Range("M134").Value = -5135.6088

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 6512392.5
$ws.Range("I4").Value = 6846246
$ws.Range("J4").Value = 2250
$ws.Range("K4").Value = 20538738
$ws.Range("L4").Value = 6750
$ws.Range("M4").Value = -20538626
$ws.Range("N4").Value = -6974
# Row 8
$ws.Range("H8").Value = 1051.5
$ws.Range("I8").Value = 1051.5
$ws.Range("K8").Value = 3154.5
$ws.Range("M8").Value = -3015.5
# Row 38
$ws.Range("H38").Value = 68.37036999999999
$ws.Range("I38").Value = 74
$ws.Range("J38").Value = 61.333332
$ws.Range("K38").Value = 222
$ws.Range("L38").Value = 183.999996
$ws.Range("M38").Value = 125
$ws.Range("N38").Value = -877.999996
# Row 104
$ws.Range("H104").Value = 6935.75
$ws.Range("I104").Value = 99
$ws.Range("J104").Value = 7912.4287
$ws.Range("K104").Value = 297
$ws.Range("L104").Value = 23737.2861
$ws.Range("M104").Value = 2324
$ws.Range("N104").Value = -28979.2861
# Row 107
$ws.Range("H107").Value = 33334564
$ws.Range("J107").Value = 1413
$ws.Range("L107").Value = 4239
$ws.Range("N107").Value = -8079
# Row 119
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").Value = ""
# Row 120
$ws.Range("H120").Value = 40000
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 40000
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 120000
$ws.Range("M120").Value = ""
$ws.Range("N120").Value = -129676
# Row 121
$ws.Range("H121").Value = 2915.95
$ws.Range("I121").Value = 366
$ws.Range("J121").Value = 3765.9333
$ws.Range("K121").Value = 1098
$ws.Range("L121").Value = 11297.7999
$ws.Range("M121").Value = 212
$ws.Range("N121").Value = -13917.7999
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = ""
# Row 137
$ws.Range("H137").Value = 1722.8572
$ws.Range("I137").Value = 1509.2307
$ws.Range("J137").Value = 4500
$ws.Range("K137").Value = 4527.6921
$ws.Range("L137").Value = 13500
$ws.Range("M137").Value = 572.3078999999998
$ws.Range("N137").Value = -23700

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 9212
$ws.Range("I70").Value = 10669.571
$ws.Range("J70").Value = 7936.625
$ws.Range("K70").Value = 10669.571
$ws.Range("L70").Value = 7936.625
$ws.Range("M70").Value = -10399.571
$ws.Range("N70").Value = -8476.625
# Row 73
$ws.Range("H73").Value = 9212
$ws.Range("I73").Value = 10669.571
$ws.Range("J73").Value = 7936.625
$ws.Range("K73").Value = 10669.571
$ws.Range("L73").Value = 7936.625
$ws.Range("M73").Value = -9733.571
$ws.Range("N73").Value = -9808.625
# Row 80
$ws.Range("H80").Value = 2132.1
$ws.Range("I80").Value = 1165.25
$ws.Range("K80").Value = 1165.25
$ws.Range("M80").Value = -167.25
# Row 83
$ws.Range("H83").Value = 2132.1
$ws.Range("I83").Value = 1165.25
$ws.Range("K83").Value = 5826.25
$ws.Range("M83").Value = -834.25
# Row 132
$ws.Range("H132").Value = 280790.12
$ws.Range("I132").Value = 305816.5
$ws.Range("K132").Value = 917449.5
$ws.Range("M132").Value = -914919.5
# Row 138
$ws.Range("H138").Value = 75499.5
$ws.Range("J138").Value = 75499.5
$ws.Range("L138").Value = 75499.5
$ws.Range("N138").Value = -85779.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 13298.154
$ws.Range("I7").Value = 13656.333
$ws.Range("K7").Value = 13656.333
$ws.Range("M7").Value = -13544.333
# Row 16
$ws.Range("H16").Value = 5972.5386
$ws.Range("I16").Value = 4464.3
$ws.Range("J16").Value = 11000
$ws.Range("K16").Value = 4464.3
$ws.Range("L16").Value = 11000
$ws.Range("M16").Value = -4294.3
$ws.Range("N16").Value = -11340
# Row 126
$ws.Range("H126").Value = 13298.154
$ws.Range("I126").Value = 13656.333
$ws.Range("K126").Value = 40968.999
$ws.Range("M126").Value = -38498.999

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 1545.3889
$ws.Range("I100").Value = 1442.909
$ws.Range("J100").Value = 1706.4286
$ws.Range("K100").Value = 2885.818
$ws.Range("L100").Value = 3412.8572
$ws.Range("M100").Value = -2344.818
$ws.Range("N100").Value = -4494.8572
# Row 107
$ws.Range("H107").Value = 1042.6666
$ws.Range("I107").Value = 807.9167
$ws.Range("K107").Value = 2423.7501
$ws.Range("M107").Value = -503.7501000000002
# Row 113
$ws.Range("H113").Value = 840.5833
$ws.Range("I113").Value = 585.2857
$ws.Range("J113").Value = 1198
$ws.Range("K113").Value = 1755.8571
$ws.Range("L113").Value = 3594
$ws.Range("M113").Value = 414.1428999999998
$ws.Range("N113").Value = -7934
# Row 126
$ws.Range("H126").Value = 3757.2273
$ws.Range("I126").Value = 3512.3333
$ws.Range("K126").Value = 10536.9999
$ws.Range("M126").Value = -8066.999899999999
# Row 136
$ws.Range("H136").Value = 2473.6875
$ws.Range("I136").Value = 2112.7856
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 6338.3568
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -3788.3568
$ws.Range("N136").Value = -20100

